# Refuse non unique subject names
# Adds a new localization entry ("enter_unique") to the "string" table on
# the "string" worksheet, inserted (alphabetically, by key) at row 40 -
# right before "french" - which pushes the existing rows 40-113 down to
# 41-114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at row 40; this shifts existing rows 40-113 down to 41-114.
$ws.Rows(40).Insert()

# Populate the new row with the new translation key + its 4 locales.
$ws.Range("A40").Value = "enter_unique"
$ws.Range("B40").Value = "Enter a unique name"
$ws.Range("C40").Value = "Entrez un nom unique"
$ws.Range("D40").Value = "Geben Sie einen eindeutigen Namen ein"
$ws.Range("E40").Value = "Gidd en eendeitegen Numm un"

# Grow the table (and its autofilter) so the new row is part of "string".
$lo.Resize($ws.Range("A1:E114")) | Out-Null

# The plain row-insert above does not shift existing hyperlink anchors, so
# fix those up by hand. Deleting any one hyperlink via a Range clears the
# whole sheet collection here, so do it once and re-add all four links at
# their correct (post-insert) locations, restoring each cell's style
# afterwards since Hyperlinks.Add stamps its own formatting on the cell.
$ws.Range("B45").Hyperlinks.Delete() | Out-Null

$r1 = $ws.Range("B46")
$s1 = $r1.Style
$ws.Hyperlinks.Add($r1, "https://github.com/NightDreamGames/Graded") | Out-Null
$r1.Style = $s1

$r2 = $ws.Range("C46:E46")
$s2 = $r2.Style
$ws.Hyperlinks.Add($r2, "https://github.com/NightDreamGames/Graded", "", "", "https://github.com/NightDreamGames/Graded") | Out-Null
$r2.Style = $s2

$r3 = $ws.Range("B38")
$s3 = $r3.Style
$ws.Hyperlinks.Add($r3, "mailto:contact@nightdreamgames.com") | Out-Null
$r3.Style = $s3

$r4 = $ws.Range("C38:E38")
$s4 = $r4.Style
$ws.Hyperlinks.Add($r4, "mailto:contact@nightdreamgames.com", "", "", "contact@nightdreamgames.com") | Out-Null
$r4.Style = $s4

# The ExternalData_1 defined name tracked the table and needs to grow too.
$wb.Names.Item(1).RefersTo = "=string!`$A`$1:`$B`$114"

# Match the saved selection/view state (whole table selected).
$ws.Range("A1:E114").Select() | Out-Null
